$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns used in every data row (constant across the whole sheet):
$colA = $ws.Cells.Item(236, 1).Value2
$colB = $ws.Cells.Item(236, 2).Value2
$colC = $ws.Cells.Item(236, 3).Value2
$colE = $ws.Cells.Item(236, 5).Value2
$colF = $ws.Cells.Item(236, 6).Value2
$colG = $ws.Cells.Item(236, 7).Value2
$colH = $ws.Cells.Item(236, 8).Value2
$colN = $ws.Cells.Item(236, 14).Value2
$colO = $ws.Cells.Item(236, 15).Value2
$colQ = $ws.Cells.Item(236, 17).Value2
$colR = $ws.Cells.Item(236, 18).Value2

# Shift the data block (rows 173..236) down by two rows, to 175..238,
# reading bottom-to-top so a source row is never clobbered before use.
for ($i = 236; $i -ge 173; $i--) {
    $dest = $i + 2

    $dVal = $ws.Cells.Item($i, 4).Value2
    $iVal = $ws.Cells.Item($i, 9).Value2
    $jVal = $ws.Cells.Item($i, 10).Value2
    $kVal = $ws.Cells.Item($i, 11).Value2
    $lVal = $ws.Cells.Item($i, 12).Value2
    $mVal = $ws.Cells.Item($i, 13).Value2
    $pVal = $ws.Cells.Item($i, 16).Value2

    $ws.Cells.Item($dest, 1).Value2 = $colA
    $ws.Cells.Item($dest, 2).Value2 = $colB
    $ws.Cells.Item($dest, 3).Value2 = $colC
    $ws.Cells.Item($dest, 4).Value2 = $dVal
    $ws.Cells.Item($dest, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($dest, 5).Value2 = $colE
    $ws.Cells.Item($dest, 6).Value2 = $colF
    $ws.Cells.Item($dest, 7).Value2 = $colG
    $ws.Cells.Item($dest, 8).Value2 = $colH
    $ws.Cells.Item($dest, 9).Value2 = $iVal
    $ws.Cells.Item($dest, 10).Value2 = $jVal
    $ws.Cells.Item($dest, 11).Value2 = $kVal
    $ws.Cells.Item($dest, 12).Value2 = $lVal
    $ws.Cells.Item($dest, 13).Value2 = $mVal
    $ws.Cells.Item($dest, 14).Value2 = $colN
    $ws.Cells.Item($dest, 15).Value2 = $colO
    $ws.Cells.Item($dest, 16).Value2 = $pVal
    $ws.Cells.Item($dest, 17).Value2 = $colQ
    $ws.Cells.Item($dest, 18).Value2 = $colR
}

# Row 173: new record (Segunda), row 174: new record (Tercera)
$ws.Cells.Item(173, 4).Value2 = 44524
$ws.Cells.Item(173, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(173, 9).Value2 = "Segunda"
$ws.Cells.Item(173, 10).Value2 = 1200
$ws.Cells.Item(173, 11).Value2 = 350
$ws.Cells.Item(173, 12).Value2 = 400
$ws.Cells.Item(173, 13).Value2 = 375
$ws.Cells.Item(173, 16).Value2 = 375

$ws.Cells.Item(174, 4).Value2 = 44524
$ws.Cells.Item(174, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(174, 9).Value2 = "Tercera"
$ws.Cells.Item(174, 10).Value2 = 1200
$ws.Cells.Item(174, 11).Value2 = 250
$ws.Cells.Item(174, 12).Value2 = 300
$ws.Cells.Item(174, 13).Value2 = 275
$ws.Cells.Item(174, 16).Value2 = 275

Write-Output "done"
